$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Test"
$ws.Range("C2").Value = "test"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "resrewr"
$ws.Range("C3").Value = "qdasdsd"

$ws.Range("C3").Select()
